# The template-migrated-expected.docx fixture splits two existing runs
# into two runs each (same text, same formatting) so the parser's new
# TokenIteratorFieldRewriterSplit can rewrite the field token piecewise:
#   "{m"                           -> "{"  + "m"
#   "'test string'.serviceA(self)}" -> "'test string'.serviceA(self)" + "}"
#
# Word has no "Run" object to split directly, so we force the run break
# the same way a user would in the UI: insert a paragraph mark at the
# split point (which always starts a new run) and then delete that mark
# again, which merges the paragraphs back together while leaving the
# run boundary in place.

$d = $word.ActiveDocument

function Split-RunAt([int]$pos) {
    $breakPoint = $d.Range($pos, $pos)
    $breakPoint.InsertParagraphAfter()
    $mark = $d.Range($pos, $pos + 1)
    $mark.Delete()
}

# Locate "{m" and split right after the "{".
$find1 = $d.Content
$find1.Find.Execute("{m", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPos1 = $find1.Start + 1
Split-RunAt $splitPos1

# Locate "(self)}" and split right before the trailing "}".
$find2 = $d.Content
$find2.Find.Execute("(self)}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPos2 = $find2.End - 1
Split-RunAt $splitPos2
